# Updates the cryptos price/volume table to the latest scrape.
# D = Price column (numeric-looking text, must stay text), E = Volume(1h) column (plain text),
# B/C rows 19-20 additionally swap (Chainlink <-> WrappedEther ranking changed).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @(
    @{ Ref = "D2"; Value = '67.063.99'; IsPrice = $true },
    @{ Ref = "E2"; Value = '  -0.35%  '; IsPrice = $false },
    @{ Ref = "D3"; Value = '3.090.10'; IsPrice = $true },
    @{ Ref = "E3"; Value = '  -1.50%  '; IsPrice = $false },
    @{ Ref = "E4"; Value = '  +0.11%  '; IsPrice = $false },
    @{ Ref = "D5"; Value = '579.12'; IsPrice = $true },
    @{ Ref = "E5"; Value = '  -0.21%  '; IsPrice = $false },
    @{ Ref = "D6"; Value = '171.24'; IsPrice = $true },
    @{ Ref = "E6"; Value = '  -2.11%  '; IsPrice = $false },
    @{ Ref = "D7"; Value = '1.00'; IsPrice = $true },
    @{ Ref = "E7"; Value = '  +0.05%  '; IsPrice = $false },
    @{ Ref = "D8"; Value = '3.086.65'; IsPrice = $true },
    @{ Ref = "E8"; Value = '  -1.48%  '; IsPrice = $false },
    @{ Ref = "D9"; Value = '0.517'; IsPrice = $true },
    @{ Ref = "E9"; Value = '  -1.50%  '; IsPrice = $false },
    @{ Ref = "D10"; Value = '6.44'; IsPrice = $true },
    @{ Ref = "E10"; Value = '  -0.97%  '; IsPrice = $false },
    @{ Ref = "D11"; Value = '0.151'; IsPrice = $true },
    @{ Ref = "E11"; Value = '  -2.92%  '; IsPrice = $false },
    @{ Ref = "D12"; Value = '0.475'; IsPrice = $true },
    @{ Ref = "E12"; Value = '  -1.87%  '; IsPrice = $false },
    @{ Ref = "D13"; Value = '0.0000244'; IsPrice = $true },
    @{ Ref = "E13"; Value = '  -2.45%  '; IsPrice = $false },
    @{ Ref = "D14"; Value = '36.44'; IsPrice = $true },
    @{ Ref = "E14"; Value = '  -2.78%  '; IsPrice = $false },
    @{ Ref = "E15"; Value = '  -1.60%  '; IsPrice = $false },
    @{ Ref = "D16"; Value = '3.604.51'; IsPrice = $true },
    @{ Ref = "E16"; Value = '  -1.38%  '; IsPrice = $false },
    @{ Ref = "D17"; Value = '67.055.88'; IsPrice = $true },
    @{ Ref = "E17"; Value = '  -0.35%  '; IsPrice = $false },
    @{ Ref = "D18"; Value = '7.06'; IsPrice = $true },
    @{ Ref = "E18"; Value = '  -1.89%  '; IsPrice = $false },
    @{ Ref = "B19"; Value = 'WrappedEther'; IsPrice = $false },
    @{ Ref = "C19"; Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'; IsPrice = $false },
    @{ Ref = "D19"; Value = '3.092.38'; IsPrice = $true },
    @{ Ref = "E19"; Value = '  -1.41%  '; IsPrice = $false },
    @{ Ref = "B20"; Value = 'Chainlink'; IsPrice = $false },
    @{ Ref = "C20"; Value = 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link'; IsPrice = $false },
    @{ Ref = "D20"; Value = '16.52'; IsPrice = $true },
    @{ Ref = "E20"; Value = '  +2.28%  '; IsPrice = $false },
    @{ Ref = "D21"; Value = '488.47'; IsPrice = $true },
    @{ Ref = "E21"; Value = '  +0.04%  '; IsPrice = $false },
    @{ Ref = "D22"; Value = '0.695'; IsPrice = $true },
    @{ Ref = "E22"; Value = '  -3.11%  '; IsPrice = $false },
    @{ Ref = "D23"; Value = '7.74'; IsPrice = $true },
    @{ Ref = "E23"; Value = '  +0.57%  '; IsPrice = $false },
    @{ Ref = "D24"; Value = '83.51'; IsPrice = $true },
    @{ Ref = "E24"; Value = '  -0.94%  '; IsPrice = $false },
    @{ Ref = "D25"; Value = '12.94'; IsPrice = $true },
    @{ Ref = "E25"; Value = '  -2.44%  '; IsPrice = $false },
    @{ Ref = "D26"; Value = '2.25'; IsPrice = $true },
    @{ Ref = "E26"; Value = '  -3.52%  '; IsPrice = $false },
    @{ Ref = "D27"; Value = '10.39'; IsPrice = $true },
    @{ Ref = "E27"; Value = '  +3.04%  '; IsPrice = $false },
    @{ Ref = "E28"; Value = '  -0.06%  '; IsPrice = $false },
    @{ Ref = "E29"; Value = '  -2.69%  '; IsPrice = $false },
    @{ Ref = "D30"; Value = '2.31'; IsPrice = $true },
    @{ Ref = "E31"; Value = '  -1.68%  '; IsPrice = $false },
    @{ Ref = "D32"; Value = '28.15'; IsPrice = $true },
    @{ Ref = "E32"; Value = '  -2.52%  '; IsPrice = $false },
    @{ Ref = "E33"; Value = '  -2.11%  '; IsPrice = $false },
    @{ Ref = "D34"; Value = '0.0₃0940'; IsPrice = $true },
    @{ Ref = "E34"; Value = '  -5.71%  '; IsPrice = $false },
    @{ Ref = "D35"; Value = '1.00'; IsPrice = $true },
    @{ Ref = "E35"; Value = '  +0.12%  '; IsPrice = $false },
    @{ Ref = "D36"; Value = '5.74'; IsPrice = $true },
    @{ Ref = "E36"; Value = '  -3.44%  '; IsPrice = $false },
    @{ Ref = "D37"; Value = '0.963'; IsPrice = $true },
    @{ Ref = "E37"; Value = '  -2.71%  '; IsPrice = $false },
    @{ Ref = "D38"; Value = '46.39'; IsPrice = $true },
    @{ Ref = "E38"; Value = '  -2.43%  '; IsPrice = $false },
    @{ Ref = "D39"; Value = '2.01'; IsPrice = $true },
    @{ Ref = "E39"; Value = '  -5.03%  '; IsPrice = $false },
    @{ Ref = "E40"; Value = '  +0.22%  '; IsPrice = $false },
    @{ Ref = "D41"; Value = '0.304'; IsPrice = $true },
    @{ Ref = "E41"; Value = '  -2.79%  '; IsPrice = $false },
    @{ Ref = "D42"; Value = '8.41'; IsPrice = $true },
    @{ Ref = "E42"; Value = '  -3.07%  '; IsPrice = $false },
    @{ Ref = "D43"; Value = '2.787.82'; IsPrice = $true },
    @{ Ref = "E43"; Value = '  -2.03%  '; IsPrice = $false },
    @{ Ref = "D44"; Value = '379.85'; IsPrice = $true },
    @{ Ref = "E44"; Value = '  -1.18%  '; IsPrice = $false },
    @{ Ref = "E45"; Value = '  -9.00%  '; IsPrice = $false },
    @{ Ref = "D46"; Value = '0.0348'; IsPrice = $true },
    @{ Ref = "E46"; Value = '  -2.92%  '; IsPrice = $false },
    @{ Ref = "D47"; Value = '135.03'; IsPrice = $true },
    @{ Ref = "E47"; Value = '  -0.66%  '; IsPrice = $false },
    @{ Ref = "D49"; Value = '24.72'; IsPrice = $true },
    @{ Ref = "E49"; Value = '  -0.98%  '; IsPrice = $false },
    @{ Ref = "D50"; Value = '2.18'; IsPrice = $true },
    @{ Ref = "E50"; Value = '  -2.10%  '; IsPrice = $false },
    @{ Ref = "E51"; Value = '  -2.39%  '; IsPrice = $false }
)

foreach ($u in $updates) {
    if ($u.IsPrice) {
        # Force text storage so number-like strings (e.g. "1.00", "67.063.99")
        # are not reinterpreted as numeric/date values, then drop the temporary
        # number-format override so no stray style survives on the cell.
        $ws.Range($u.Ref).NumberFormat = "@"
        $ws.Range($u.Ref).Value = $u.Value
        $ws.Range($u.Ref).Style = "Normal"
    } else {
        $ws.Range($u.Ref).Value = $u.Value
    }
}
